$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Hunk 1 + Hunk 3: "{{sigfirm|req|lawyer 01}}" and "{{sigfirm|req|lawyer 02...}}"
# paragraphs still have the role name split across three runs
# (<w:t>sig</w:t> / <w:t>firm</w:t> / <w:t>|req|lawyer</w:t>). Re-saving the
# exact same text over that span merges the three runs back into one, which
# is exactly what the target markup needs (and leaves the paragraphs whose
# role name is already a single run untouched).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("sigfirm|req|lawyer", $true, $false, $false, $false, $false, $true, 1, $false, "sigfirm|req|lawyer", 2) | Out-Null

# ---------------------------------------------------------------------------
# Hunk 2: the "{{sig|req|lawyer 02            }}" paragraph (no "firm" in the
# role name) keeps its trailing " 02            }}" as a single run in the
# source, but the target splits it into " " / "test" / "            }}".
# Locate that paragraph uniquely via its text, then clear and retype the
# span in three pieces so Word creates three separate runs.
# ---------------------------------------------------------------------------
$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $para = $d.Paragraphs($i)
    $t = $para.Range.Text
    if ($t.StartsWith("{{") -and $t -like "*sig|req|lawyer 02*") {
        $r = $para.Range
        $pStart = $r.Start
        $full = $r.Text
        $relIdx = $full.IndexOf(" 02            }}")
        $absStart = $pStart + $relIdx
        $runRange = $d.Range($absStart, $absStart + (" 02            }}".Length))
        $runRange.Text = ""
        $ip1 = $d.Range($absStart, $absStart)
        $ip1.InsertAfter(" ")
        $absStart2 = $absStart + 1
        $ip2 = $d.Range($absStart2, $absStart2)
        $ip2.InsertAfter("test")
        $absStart3 = $absStart2 + 4
        $ip3 = $d.Range($absStart3, $absStart3)
        $ip3.InsertAfter("            }}")
        break
    }
}

# ---------------------------------------------------------------------------
# Hunk 4: the "{{sigfirm|req|lawyer 03            }}" paragraph has its
# trailing text split as " 0" / "3" / "            }}" in the source; the
# target collapses it back into a single " 03            }}" run.
# ---------------------------------------------------------------------------
$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $para = $d.Paragraphs($i)
    $t = $para.Range.Text
    if ($t.StartsWith("{{") -and $t -like "*sigfirm|req|lawyer 03*") {
        $r = $para.Range
        $pStart = $r.Start
        $full = $r.Text
        $relIdx = $full.IndexOf(" 03            }}")
        $absStart = $pStart + $relIdx
        $runRange = $d.Range($absStart, $absStart + (" 03            }}".Length))
        $runRange.Text = ""
        $ip = $d.Range($absStart, $absStart)
        $ip.InsertAfter(" 03            }}")
        break
    }
}

Write-Output "done"
